# إضافة حدث جديد في Card24
# The "Card24" sheet keeps a history of service events. The most recently
# added row (16) was still an empty placeholder (only the date/event/
# correction/serviced-by columns were filled in). This edit:
#   1. Finalises row 16 by stamping its still-blank numeric columns
#      (B..K and P) with the literal placeholder text "nan", matching every
#      other already-finalised row above it.
#   2. Appends a brand-new empty placeholder row (17) underneath, copied
#      from row 16 so it keeps the same card/date/event/correction/
#      serviced-by values and cell formatting, ready for the next event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$lastRow = 16
$newRow = $lastRow + 1

# Duplicate the whole row (values + formatting) down one row first, so the
# brand-new row inherits row 16's current content/format untouched -
# including the still-empty B:K / P placeholder cells.
$ws.Range("A" + $lastRow + ":P" + $lastRow).Copy($ws.Range("A" + $newRow + ":P" + $newRow))

# Now finalise row 16: its previously-empty numeric/observation columns get
# the "nan" placeholder text used throughout the rest of the sheet.
$ws.Range("B" + $lastRow + ":K" + $lastRow).Value = "nan"
$ws.Cells.Item($lastRow, 16).Value = "nan"
